# Auto-generated edit script: update cryptos list values (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "25.995.43"
$ws.Cells.Item(2, 5).Value = "  -0.75%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.646.27"
$ws.Cells.Item(3, 5).Value = "  -0.68%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.26%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'217.63"
$ws.Cells.Item(5, 5).Value = "  -0.57%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'0.5226"
$ws.Cells.Item(6, 5).Value = "  +0.42%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.19%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.2616"
$ws.Cells.Item(8, 5).Value = "  -1.94%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.06265"
$ws.Cells.Item(9, 5).Value = "  -0.94%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'20.45"
$ws.Cells.Item(10, 5).Value = "  -3.16%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.07742"
$ws.Cells.Item(11, 5).Value = "  -0.16%  "

# Row 12
$ws.Cells.Item(12, 2).Value = "Polkadot"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(12, 4).Value = "'4.456"
$ws.Cells.Item(12, 5).Value = "  +0.42%  "

# Row 13
$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13, 4).Value = "1.650.08"
$ws.Cells.Item(13, 5).Value = "  -0.56%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'0.5434"
$ws.Cells.Item(14, 5).Value = "  -0.76%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "0.0₅8062"
$ws.Cells.Item(15, 5).Value = "  -2.17%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'64.66"
$ws.Cells.Item(16, 5).Value = "  -0.49%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "26.004.11"
$ws.Cells.Item(17, 5).Value = "  -0.86%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  -0.23%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'4.546"
$ws.Cells.Item(19, 5).Value = "  -2.84%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'191.15"
$ws.Cells.Item(20, 5).Value = "  -0.49%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  -1.44%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'5.973"
$ws.Cells.Item(22, 5).Value = "  -2.35%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -0.30%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'139.01"
$ws.Cells.Item(24, 5).Value = "  +1.15%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'0.1232"
$ws.Cells.Item(25, 5).Value = "  -0.28%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'7.246"
$ws.Cells.Item(26, 5).Value = "  -0.24%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'16.11"
$ws.Cells.Item(27, 5).Value = "  -0.05%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'1.427"
$ws.Cells.Item(28, 5).Value = "  +1.03%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'0.05915"
$ws.Cells.Item(29, 5).Value = "  -1.88%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'1.272"
$ws.Cells.Item(30, 5).Value = "  -1.18%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'3.491"
$ws.Cells.Item(31, 5).Value = "  -1.69%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'3.231"
$ws.Cells.Item(32, 5).Value = "  -3.55%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'1.515"
$ws.Cells.Item(33, 5).Value = "  -8.10%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'2.416"
$ws.Cells.Item(34, 5).Value = "  +0.14%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'0.9422"
$ws.Cells.Item(35, 5).Value = "  -4.13%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "'2.749"
$ws.Cells.Item(36, 5).Value = "  -1.06%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'0.5656"
$ws.Cells.Item(37, 5).Value = "  -4.83%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'0.01602"
$ws.Cells.Item(38, 5).Value = "  +0.66%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'5.843"
$ws.Cells.Item(39, 5).Value = "  -2.16%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'0.8452"

# Row 41
$ws.Cells.Item(41, 5).Value = "  -0.14%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'100.45"
$ws.Cells.Item(42, 5).Value = "  +0.62%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "1.000.32"
$ws.Cells.Item(43, 5).Value = "  -3.82%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "1.789.13"
$ws.Cells.Item(44, 5).Value = "  -0.53%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(45, 4).Value = "0.0₈107"
$ws.Cells.Item(45, 5).Value = "  -1.93%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "Aave"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(46, 4).Value = "'56.53"
$ws.Cells.Item(46, 5).Value = "  -1.14%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -0.40%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +1.52%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "RenderToken"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(49, 4).Value = "'1.475"
$ws.Cells.Item(49, 5).Value = "  +0.05%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'0.05146"
$ws.Cells.Item(50, 5).Value = "  -0.68%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "EnergySwap"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(51, 4).Value = "'7.821"
$ws.Cells.Item(51, 5).Value = "  -3.64%  "
